# add driver name on laporan detail -> kenek sheet
$wb = $excel.ActiveWorkbook

# "komisi kenek" sheet: insert a new "Driver" column between "Kenek" (C) and
# "No. Polisi" (old D, now E), shifting the remaining columns right.
$ws = $wb.Worksheets.Item("komisi kenek")
$ws.Columns("D:D").Insert() | Out-Null
$ws.Range("D1").Value = "Driver"

# Keep the (hidden) filter-database defined names for "komisi kenek" in sync
# with the new last column (I -> J).
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "komisi kenek!*") {
        $n.RefersTo = "='komisi kenek'!`$A`$1:`$J`$1"
    }
}

# Make "komisi kenek" the active sheet/tab and leave the selection on A2,
# matching the latest saved view state.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
